$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SampleServiceNew_1_0")
$ws.Select()
$ws.Range("D2").Value = "<?xml version=`"1.0`" encoding=`"UTF-8`"?>`n<FailureResponse>MSISDN is not proper</FailureResponse>"
$ws.Range("B11").Select()
